$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "OC140 sec"

$ws.Range("H2").Value = 9.1300000000000008
$ws.Range("H3").Value = 0.1
$ws.Range("H4").Value = 0.53
$ws.Range("H5").Value = 0.96
$ws.Range("H6").Value = 2160
$ws.Range("H7").Value = 2022000

$ws.Range("H8").Formula = "=1000000*H4/H6"
$ws.Range("H9").Formula = "=(H5-H4)/(H10*H6)"
$ws.Range("H10").Formula = "=(H2-H3)/H7"
$ws.Range("H11").Formula = "=H10*1000000"

$ws.Range("H8").NumberFormat = $ws.Range("F8").NumberFormat
$ws.Range("H9").NumberFormat = $ws.Range("F9").NumberFormat
$ws.Range("H11").NumberFormat = $ws.Range("F11").NumberFormat

$ws.Range("I8").NumberFormat = $ws.Range("F8").NumberFormat
$ws.Range("I9").NumberFormat = $ws.Range("F9").NumberFormat
$ws.Range("I11").NumberFormat = $ws.Range("F11").NumberFormat

$null = $ws.Range("H9").Select()
